# Refresh the "cryptos" price table (Price / Volume(1h) columns, plus the
# Cronos/HuobiToken row swap) to match the latest scrape, as produced by the
# "Updated cryptos list ... with GitHub Actions" job.
#
# NOTE: every "Price" cell in column D is stored as literal TEXT (some values
# use dotted thousand separators like "37.050.18" that aren't valid numbers
# anyway). Values that *do* look like plain numbers (e.g. "227.26") are
# written with a leading apostrophe so Excel keeps them as text instead of
# silently coercing them to a Number (which would also drop significant
# trailing zeros, e.g. "14.20" -> 14.2).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '37.105.22'
$ws.Range('E2').Value = '  -0.03%  '
$ws.Range('D3').Value = '2.015.95'
$ws.Range('E3').Value = '  -0.32%  '
$ws.Range('E4').Value = '  -0.14%  '
$ws.Range('D5').Value = '''227.26'
$ws.Range('E5').Value = '  +0.02%  '
$ws.Range('E6').Value = '  +0.32%  '
$ws.Range('E7').Value = '  -0.05%  '
$ws.Range('D8').Value = '''55.73'
$ws.Range('E8').Value = '  +1.29%  '
$ws.Range('E9').Value = '  -1.55%  '
$ws.Range('D10').Value = '''0.0779'
$ws.Range('E10').Value = '  -1.26%  '
$ws.Range('E11').Value = '  -2.71%  '
$ws.Range('D12').Value = '2.314.56'
$ws.Range('D13').Value = '''14.20'
$ws.Range('D14').Value = '''19.92'
$ws.Range('E14').Value = '  -2.31%  '
$ws.Range('E15').Value = '  -1.08%  '
$ws.Range('D16').Value = '''5.16'
$ws.Range('E16').Value = '  -0.02%  '
$ws.Range('D17').Value = '2.015.69'
$ws.Range('E17').Value = '  -0.46%  '
$ws.Range('D18').Value = '37.012.28'
$ws.Range('E18').Value = '  -0.03%  '
$ws.Range('E19').Value = '  +1.70%  '
$ws.Range('E20').Value = '  +0.07%  '
$ws.Range('D21').Value = '0.0₃0814'
$ws.Range('E21').Value = '  -2.21%  '
$ws.Range('D22').Value = '''222.51'
$ws.Range('E22').Value = '  -0.20%  '
$ws.Range('E23').Value = '  -0.22%  '
$ws.Range('E24').Value = '  +1.36%  '
$ws.Range('E25').Value = '  -3.46%  '
$ws.Range('D26').Value = '''163.23'
$ws.Range('E26').Value = '  -2.67%  '
$ws.Range('E27').Value = '  -4.07%  '
$ws.Range('E28').Value = '  +1.77%  '
$ws.Range('D29').Value = '''18.67'
$ws.Range('E29').Value = '  -0.18%  '
$ws.Range('E30').Value = '  -1.90%  '
$ws.Range('E31').Value = '  +0.28%  '
$ws.Range('D32').Value = '''4.43'
$ws.Range('E32').Value = '  -1.02%  '
$ws.Range('E33').Value = '  -0.67%  '
$ws.Range('E34').Value = '  -0.09%  '
$ws.Range('E35').Value = '  -0.53%  '
$ws.Range('E36').Value = '  +2.67%  '
$ws.Range('E37').Value = '  -0.18%  '
$ws.Range('D38').Value = '''3.16'
$ws.Range('E38').Value = '  +0.41%  '
$ws.Range('D39').Value = '''5.45'
$ws.Range('E39').Value = '  +1.76%  '
$ws.Range('D40').Value = '1.468.12'
$ws.Range('E40').Value = '  -1.67%  '
$ws.Range('E41').Value = '  -2.57%  '
$ws.Range('D42').Value = '''4.29'
$ws.Range('E42').Value = '  +17.04%  '
$ws.Range('D43').Value = '''93.85'
$ws.Range('E43').Value = '  -0.96%  '
$ws.Range('B44').Value = 'HuobiToken'
$ws.Range('C44').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$ws.Range('D44').Value = '''2.79'
$ws.Range('E44').Value = '  -2.24%  '
$ws.Range('B45').Value = 'Cronos'
$ws.Range('C45').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D45').Value = '''0.0910'
$ws.Range('E45').Value = '  -2.13%  '
$ws.Range('D46').Value = '''16.20'
$ws.Range('E46').Value = '  -1.80%  '
$ws.Range('E47').Value = '  -2.21%  '
$ws.Range('E48').Value = '  +0.01%  '
$ws.Range('E49').Value = '  +0.39%  '
$ws.Range('E50').Value = '  -1.17%  '
$ws.Range('D51').Value = '2.203.80'
$ws.Range('E51').Value = '  -0.48%  '
